# Actualizacion automatica del mapa (2025-10-30 11:51:58)
#
# Inserts one new incident row before the current row 83 (shifting the
# existing rows 83-88 down to 84-89) and then overwrites the final row
# (89) — which used to be the old row 88 — with a brand-new record, since
# that old trailing record is being dropped from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing TEXT storage (so numeric
# looking strings like "-657", "13", "810454540" or date-like strings such
# as "10/27/2025" are kept as literal text instead of being auto-converted
# to a number/date by Excel), and leave the cell's style unchanged
# afterwards (matches the source workbook, where data rows carry no
# explicit style).
function Set-TextCell($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Helper: write a row of record data (matches the 18 columns A:R used by
# this sheet) starting at the given row number.
function Set-RowData($rowNum, $values) {
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = $values[$i]
        $addr = "$col$rowNum"
        if ($col -eq "I" -or $col -eq "M" -or $col -eq "N") {
            # Numeric columns: Attachments (I), Coordenada_X (M), Coordenada_Y (N)
            $ws.Range($addr).Value = $val
        } else {
            Set-TextCell $ws.Range($addr) $val
        }
    }
}

# 1) Insert a brand-new row above row 83; everything from 83 downward
#    shifts down by one (old 83->84, 84->85, ..., 88->89).
$ws.Rows.Item(83).Insert()

# 2) Populate the freshly inserted row 83 with its data.
Set-RowData 83 @(
    "-657",
    "10/27/2025",
    "Conde 1632",
    "13",
    "810454540",
    "Optical Power",
    "Pendiente",
    "Poste inclinado cambiar o desmontar",
    1,
    "Cambio",
    "Sin equipos",
    "Poste",
    -58.461492,
    -34.57199,
    "Colegiales",
    "Capital Norte",
    "ATH-R",
    "Fuera de Poligono OVL"
)

# 3) The old trailing row (previously row 88, now shifted to row 89) is
#    replaced outright with a new record -- clear it first, then write
#    the new values.
$ws.Rows.Item(89).ClearContents()
Set-RowData 89 @(
    "-661",
    "10/30/2025",
    "FLORES, VENANCIO, GRAL. 3715",
    "10",
    "",
    "Optical Power",
    "Pendiente",
    "Picada",
    1,
    "Cambio",
    "Sin equipos",
    "Pasante",
    -58.47989,
    -34.631709,
    "Devoto",
    "Capital Norte",
    "DEV-L",
    "ARATO-25058.PO.2DEV"
)
